# Daily attendance processing - 2025-11-15 17:18:57
#
# For every data row in the "Recorded By" column (G), the list of
# recorder names/emails (comma-separated) is reversed in order.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# e.g. "system, backup@backdoor.com, System" -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
